$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap a few country names that moved position in the source data ---
# (Eslovaquia / Eslovenia)
$ws.Range("A89").Value = "Eslovenia"
$ws.Range("A90").Value = "Eslovaquia"

# (Belice / Nueva Caledonia)
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("A193").Value = "Belice"

# (Butan / Islas Virgenes Britanicas)
$ws.Range("A212").Value = "Islas Virgenes Britanicas"
$ws.Range("A213").Value = "Butan"

# --- Updated COVID figures (countries & provincias Spain refresh) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1348485
$ws.Range("C4").Value = 1176
$ws.Range("D4").Value = 238081
$ws.Range("E4").Value = 1030336
$ws.Range("G4").Value = 31
$ws.Range("H4").Value = 80068

# Row 5 - España
$ws.Range("F5").Value = 1650

# Row 75 - Uzbekistan
$ws.Range("D75").Value = 1881
$ws.Range("E75").Value = 520

# Row 89 - Eslovenia (after swap)
$ws.Range("C89").Value = 3
$ws.Range("D89").Value = 256
$ws.Range("E89").Value = 1099
$ws.Range("F89").Value = 10
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 102

# Row 90 - Eslovaquia (after swap)
$ws.Range("B90").Value = 1457
$ws.Range("C90").Value = 2
$ws.Range("D90").Value = 941
$ws.Range("E90").Value = 490
$ws.Range("F90").Value = 5
$ws.Range("H90").Value = 26

# Row 129 - Estado de Palestina
$ws.Range("D129").Value = 263
$ws.Range("E129").Value = 110

# Row 192 - Nueva Caledonia (after swap)
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# Row 193 - Belice (after swap)
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Row 212 - Islas Virgenes Britanicas (after swap)
$ws.Range("D212").Value = 4
$ws.Range("H212").Value = 1

# Row 213 - Butan (after swap)
$ws.Range("D213").Value = 5
$ws.Range("H213").Value = 0
